# Insert a new, empty bookmark named "bookmark_question_11_non" into the
# empty (underlined/bold) paragraph that immediately precedes the
# "DEUXIEME RESOLUTION" heading. Word renumbers every later w:bookmarkStart /
# w:bookmarkEnd id by +1 automatically when the package is saved, which is
# exactly the id shift (49->50, 50->51, ... 78->79) shown in the diff.

$d = $word.ActiveDocument

# Locate "DEUXIEME RESOLUTION" - this leaves a live Range positioned on it.
$rng = $d.Content
$found = $rng.Find.Execute("DEUXIEME RESOLUTION", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The paragraph that holds the heading text.
    $headingPara = $rng.Paragraphs(1)

    # The paragraph right before it is the empty, underlined paragraph
    # where the new bookmark belongs.
    $targetPara = $headingPara.Previous(1)

    # Add a collapsed (empty) bookmark spanning that paragraph's mark,
    # matching the existing "empty" bookmarks in the document (e.g.
    # autre_per_2, bookmark_question_10_1).
    $d.Bookmarks.Add("bookmark_question_11_non", $targetPara.Range)
}
